$d = $word.ActiveDocument

# Collapse to the very end of the document body content (right before
# sectPr) — this is where the new "Ejercicio 11" section is appended,
# immediately after the paragraph ending in "...middleware (thunk). ".
$r = $d.Content
$r.Collapse(0)

# Build the three new paragraphs (blank spacer, bold "Ejercicio 11"
# heading, and the justified body paragraph) as a raw OOXML fragment so
# the run/paragraph-mark formatting matches exactly, then inject it.
$xmlFragment = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ejercicio 11</w:t></w:r><w:r><w:t xml:space="preserve"> (3h aprox)</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>En este ejercicio he creado un modal nuevo que se emplea como formulario. El mayor reto ha sido ajustar los estilos del modal para que se vea como</w:t></w:r><w:r><w:t xml:space="preserve"> en</w:t></w:r><w:r><w:t xml:space="preserve"> el vídeo. Una vez gestionada la apertura del modal mediante una variable de estado, he gestionado en redux la inserción del comentario de la misma manera que lo hice para el marcado de un favorito, agregando los tipos de acciones pertinentes y un nuevo caso en el reducer Comentarios, además de </w:t></w:r><w:r><w:t>las funciones en el ActionCreators</w:t></w:r><w:r><w:t xml:space="preserve">. De esta manera, al rellenar el formulario y darle a enviar se agrega un nuevo comentario a la base de datos y se muestra en la aplicación. He agregado una validación en el formulario para que no deje enviar si no se han rellenado los campos de autor y comentario. </w:t></w:r></w:p>'
$r.InsertXML($xmlFragment)

Write-Output "Ejercicio 11 section inserted"
